# Applies the changes described by the commit diff to the analysis workbook.
# The workbook has two sheets: OP-0 and OP-1.

$wb = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item("OP-0")
$ws1 = $wb.Worksheets.Item("OP-1")

# ----------------------------------------------------------------------
# Sheet OP-0
# ----------------------------------------------------------------------

# 1. C30: collapse duplicated moveq/movne lines
$ws0.Range("C30").Value2 = "cmp r3, 0`nmoveq r3, 1`nmovne r3, 0"

# 2-4. Rows 41-43: remove the extraneous K/L/N "Lines Missed" entries
$ws0.Range("K41:N41").ClearContents()
$ws0.Range("K42:N42").ClearContents()
$ws0.Range("K43:N43").ClearContents()

# 5. Row 44: remove K/L/N entries and fix the line-number label
$ws0.Range("K44:N44").ClearContents()
$ws0.Range("B44").Value2 = "28, 30"

# 6-8. Rows 45-47: fix line-number labels
$ws0.Range("B45").Value2 = "28, 31"
$ws0.Range("B46").Value2 = "28, 32"
$ws0.Range("B47").Value2 = "28, 33"

# 9. C53: collapse duplicated moveq/movne lines (same pattern as C30)
$ws0.Range("C53").Value2 = "cmp r3, 0`nmoveq r3, 1`nmovne r3, 0"

# 10-12. Rows 99-101: remove the extraneous K/L/N "Lines Missed" entries
$ws0.Range("K99:N99").ClearContents()
$ws0.Range("K100:N100").ClearContents()
$ws0.Range("K101:N101").ClearContents()

# 13. Row 102: remove K/L/N entries and fix the line-number label
$ws0.Range("K102:N102").ClearContents()
$ws0.Range("B102").Value2 = "30, 32"

# 14-16. Rows 103-105: fix line-number labels
$ws0.Range("B103").Value2 = "30, 33"
$ws0.Range("B104").Value2 = "30, 34"
$ws0.Range("B105").Value2 = "30, 35"

# 17. Row 136: remove the extraneous K/L/M "Lines Missed" entries
$ws0.Range("K136:M136").ClearContents()

# 18. Row 137: the True Positive flag was mis-classified under
#     "False Positive" (H); move it to the correct "True Positives" column (G)
$ws0.Range("G137").Value = $true
$ws0.Range("H137").ClearContents()

# ----------------------------------------------------------------------
# Sheet OP-1
# ----------------------------------------------------------------------

# 19-20. Rows 33-34: same False-Positive -> True-Positive column fix
$ws1.Range("G33").Value = $true
$ws1.Range("H33").ClearContents()

$ws1.Range("G34").Value = $true
$ws1.Range("H34").ClearContents()

# 21. C69: collapse duplicated moveq line
$ws1.Range("C69").Value2 = "cmp r0, 0`nmoveq r0, 1`nbeq .L1"

# 22. C70: collapse duplicated movgt/movle lines
$ws1.Range("C70").Value2 = "cmp r0, 3`nmovgt r0, 0`nmovle r0, 1"

# 23. C76: collapse duplicated movne line
$ws1.Range("C76").Value2 = "cmp r0, 1`nmovne r0, 0`nbne .L16"

# 24. B77: expand the line-number range label
$ws1.Range("B77").Value2 = "75, 76, 78"
